$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row (row 11): right-answer marks 5 -> 4, wrong-answer penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update the "Total" row (row 12): recalculated totals based on new marking scheme
$ws.Range("B12").Value = 108
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "106 / 112"
